$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn F7:F15 into incrementing formulas based on the cell above, mirroring
# the "LOAD STEP" column building up from the literal value in F6.
$ws.Range("F7").Formula = "=F6+1"
$ws.Range("F8:F15").Formula = "=F7+1"

# Update the selection to match the final cursor position left by the author.
$ws.Range("M14").Select()
